# TC22_Canine_Filter_Breed-Giant.xlsx
# Add a new "StatQuery" column (new column B) holding the Neo4j stat-bar
# query, pushing the existing dbExcel/WebExcel columns one slot to the
# right (old B -> C, old C -> D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at B; this shifts the old B (dbExcel) and
#     old C (WebExcel) columns one place to the right. ---
[void]$ws.Columns.Item(2).Insert()

# --- Header for the new column ---
$ws.Range("B1").Value = "StatQuery"

# --- New stat-bar Neo4j query text for the new column's data row ---
$statQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Giant Schnauzer'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

$ws.Range("B2").Value = $statQuery

# Match the wrapped-text formatting used by the other long-query cell (A2)
$ws.Range("B2").WrapText = $true

# New column should be the same width as column A (75.81640625 chars)
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# --- Update the view: scroll so row 2 is at the top and select B2 ---
[void]$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
